$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1173
$ws.Range("F6").Value = 0
$ws.Range("F9").Value = 189
$ws.Range("F10").Value = 325
$ws.Range("F13").Value = 6441
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 438
$ws.Range("F17").Value = 618
$ws.Range("F19").Value = 0
$ws.Range("F23").Value = 10381
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 2271
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 2231
$ws.Range("F30").Value = 86
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 2144
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 64
$ws.Range("F38").Value = 5353
$ws.Range("F40").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("F45").Value = 1080
$ws.Range("F46").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 0

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F7").Value = 240
$ws.Range("F9").Value = 47
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F15").Value = 100
$ws.Range("F16").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F20").Value = 7

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 1173
$ws.Range("F5").Value = 9217
$ws.Range("F6").Value = 7190
$ws.Range("F7").Value = 189
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 6441
$ws.Range("F15").Value = 6441
$ws.Range("F16").Value = 1106
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 196
$ws.Range("F27").Value = 10381
$ws.Range("F28").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 2231
$ws.Range("F32").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 1431
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 5353
$ws.Range("F41").Value = 436
$ws.Range("F42").Value = 1212
$ws.Range("F43").Value = 710
$ws.Range("F44").Value = 0
$ws.Range("F47").Value = 1080
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 1399
$ws.Range("F50").Value = 0
$ws.Range("F51").Value = 1102
